$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.34819999999999
$ws.Range("B7").Value = 4.542199999999999
$ws.Range("A8").Value = -22.38390000000002
$ws.Range("A10").Value = -21.66729999999999
$ws.Range("A12").Value = -21.53280000000001
$ws.Range("B15").Value = 4.516299999999996
$ws.Range("A18").Value = -21.72900000000001
$ws.Range("B18").Value = 6.486099999999994
$ws.Range("E18").Value = 17.92420000000002
$ws.Range("E19").Value = 16.3964
$ws.Range("B20").Value = 9.009099999999998
$ws.Range("E27").Value = 16.40719999999999
$ws.Range("B29").Value = 5.178900000000003
$ws.Range("B30").Value = 4.9165
$ws.Range("B31").Value = 4.514299999999999
$ws.Range("E31").Value = 16.7223
$ws.Range("A37").Value = -20.01179999999999
$ws.Range("E38").Value = 15.99709999999999
$ws.Range("B40").Value = 9.180599999999993
$ws.Range("E42").Value = 16.4385
$ws.Range("E44").Value = 16.2366
$ws.Range("E47").Value = 16.38619999999999
$ws.Range("B50").Value = 4.633700000000001
$ws.Range("A55").Value = -21.79509999999999
$ws.Range("E58").Value = 16.34600000000002
$ws.Range("E65").Value = 17.33560000000001
$ws.Range("A68").Value = -21.5278
$ws.Range("B68").Value = 4.499800000000002
$ws.Range("E73").Value = 17.39580000000001
$ws.Range("B76").Value = 5.579599999999998
$ws.Range("A77").Value = -20.55689999999999
$ws.Range("A78").Value = -20.03759999999998
$ws.Range("A81").Value = -21.85530000000001
$ws.Range("A82").Value = -21.78800000000001
$ws.Range("B87").Value = 4.645899999999994
$ws.Range("B88").Value = 4.349199999999996
$ws.Range("E90").Value = 16.17279999999998
$ws.Range("E94").Value = 18.89570000000003
$ws.Range("E95").Value = 18.29470000000001
$ws.Range("B96").Value = 5.417300000000004
$ws.Range("B98").Value = 5.713299999999999
$ws.Range("B101").Value = 9.495299999999995
$ws.Range("E101").Value = 16.71560000000002
$ws.Range("B102").Value = 8.519400000000005